$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H1").Value = "cluster_class"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New cluster_class values for rows 2-63
$values = @(0,0,-1,0,1,1,0,1,1,2,1,1,1,2,0,-1,2,3,1,-1,1,0,1,0,1,1,2,1,2,0,1,0,3,1,1,2,1,0,-1,1,2,1,1,1,1,3,1,1,2,1,1,2,1,0,1,3,1,0,3,-1,2,0.97)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
